# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff run:
#   - new report UUID: 778d4182-e38d-4db5-b831-0f4a7e156d10 -> 60bc0b94-27e8-4675-b7fc-e32f15124f83
#   - new xlf content hash: 617723a5ac8626ca4f97dd41b8746109f019c48d -> e5743c34d1759be6b0bb167abfa7558e469bf267
#   - refreshed handoff timestamps

$wb = $excel.ActiveWorkbook

$oldGuid = "778d4182-e38d-4db5-b831-0f4a7e156d10"
$newGuid = "60bc0b94-27e8-4675-b7fc-e32f15124f83"

$oldHash = "617723a5ac8626ca4f97dd41b8746109f019c48d"
$newHash = "e5743c34d1759be6b0bb167abfa7558e469bf267"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldZh = "$oldGuid.$oldHash.zh-cn.xlf"
$newZh = "$newGuid.$newHash.zh-cn.xlf"

$oldDe = "$oldGuid.$oldHash.de-de.xlf"
$newDe = "$newGuid.$newHash.de-de.xlf"

$oldOverviewDate = "2016-03-21 22:56:42"
$newOverviewDate = "2016-03-21 22:57:09"

$oldZhDate = "2016-03-21 22:56:38"
$newZhDate = "2016-03-21 22:57:06"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsOverview.Range("D2").Value = $newOverviewDate

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsZh.Range("D2").Value = $newZh
$wsZh.Hyperlinks.Item(2).TextToDisplay = $newZh
$wsZh.Range("E2").Value = $newZhDate

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsDe.Range("D2").Value = $newDe
$wsDe.Hyperlinks.Item(2).TextToDisplay = $newDe
$wsDe.Range("E2").Value = $newOverviewDate
